$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "latitude"/"longitude" counters block had duplicated rows: delete the
# extra duplicate rows 24-27 (2x "latitude" + 2x "longitude") so the
# remaining rows below shift up by 4 and the sheet ends at row 44 instead
# of 48.
$ws.Range("A24:C27").EntireRow.Delete()
